$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 1820
$ws.Range("E2").Value = 178
$ws.Range("F2").Value = 178
$ws.Range("G2").Value = 211
$ws.Range("H2").Value = 156
$ws.Range("I2").Value = 157
$ws.Range("J2").Value = -1
$ws.Range("K2").Value = 1806
$ws.Range("L2").Value = 461
$ws.Range("M2").Value = 1345
$ws.Range("N2").Value = 1339
$ws.Range("O2").Value = 6
$ws.Range("P2").Value = 82
$ws.Range("Q2").Value = 174
$ws.Range("R2").Value = -139
$ws.Range("S2").Value = -13
$ws.Range("T2").Value = 158
$ws.Range("U2").Value = 16
$ws.Range("V2").Value = 30
$ws.Range("W2").Value = 9.779999999999999
$ws.Range("X2").Value = 8.57
$ws.Range("Y2").Value = 12.27
$ws.Range("Z2").Value = 9.16
$ws.Range("AA2").Value = 34.28
$ws.Range("AB2").Value = 1543.21
$ws.Range("AC2").Value = 800
$ws.Range("AD2").Value = 15.24
$ws.Range("AE2").Value = 6864
$ws.Range("AF2").Value = 1.78
$ws.Range("AG2").Value = 125
$ws.Range("AH2").Value = 1.02
$ws.Range("AJ2").Value = 19662677
$ws.Range("AI2").ClearContents()

# Row 3
$ws.Range("D3").Value = 2162
$ws.Range("E3").Value = 231
$ws.Range("F3").Value = 231
$ws.Range("G3").Value = 232
$ws.Range("H3").Value = 173
$ws.Range("I3").Value = 177
$ws.Range("J3").Value = -4
$ws.Range("K3").Value = 1947
$ws.Range("L3").Value = 461
$ws.Range("M3").Value = 1487
$ws.Range("N3").Value = 1482
$ws.Range("O3").Value = 5
$ws.Range("P3").Value = 86
$ws.Range("Q3").Value = 258
$ws.Range("R3").Value = -227
$ws.Range("S3").Value = -53
$ws.Range("T3").Value = 198
$ws.Range("U3").Value = 60
$ws.Range("W3").Value = 10.68
$ws.Range("X3").Value = 8.02
$ws.Range("Y3").Value = 12.55
$ws.Range("Z3").Value = 9.24
$ws.Range("AA3").Value = 31
$ws.Range("AB3").Value = 1632.26
$ws.Range("AC3").Value = 900
$ws.Range("AD3").Value = 17.18
$ws.Range("AE3").Value = 7596
$ws.Range("AF3").Value = 2.04
$ws.Range("AG3").Value = 175
$ws.Range("AH3").Value = 1.13
$ws.Range("AJ3").Value = 19662677
$ws.Range("V3").ClearContents()
$ws.Range("AI3").ClearContents()

# Row 4
$ws.Range("D4").Value = 2407
$ws.Range("E4").Value = 291
$ws.Range("F4").Value = 291
$ws.Range("G4").Value = 278
$ws.Range("H4").Value = 197
$ws.Range("I4").Value = 202
$ws.Range("J4").Value = -5
$ws.Range("K4").Value = 2147
$ws.Range("L4").Value = 527
$ws.Range("M4").Value = 1620
$ws.Range("N4").Value = 1614
$ws.Range("O4").Value = 6
$ws.Range("P4").Value = 88
$ws.Range("Q4").Value = 208
$ws.Range("R4").Value = -202
$ws.Range("S4").Value = -58
$ws.Range("T4").Value = 173
$ws.Range("U4").Value = 35
$ws.Range("W4").Value = 12.11
$ws.Range("X4").Value = 8.19
$ws.Range("Y4").Value = 13.08
$ws.Range("Z4").Value = 9.630000000000001
$ws.Range("AA4").Value = 32.51
$ws.Range("AB4").Value = 1765.73
$ws.Range("AC4").Value = 1029
$ws.Range("AD4").Value = 17.15
$ws.Range("AE4").Value = 8332
$ws.Range("AF4").Value = 2.12
$ws.Range("AG4").Value = 234
$ws.Range("AH4").Value = 1.32
$ws.Range("AJ4").Value = 19662677
$ws.Range("V4").ClearContents()
$ws.Range("AI4").ClearContents()

# Row 5
$ws.Range("D5").Value = 2655
$ws.Range("E5").Value = 253
$ws.Range("F5").Value = 253
$ws.Range("G5").Value = 235
$ws.Range("H5").Value = 103
$ws.Range("I5").Value = 107
$ws.Range("J5").Value = -5
$ws.Range("K5").Value = 2272
$ws.Range("L5").Value = 591
$ws.Range("M5").Value = 1681
$ws.Range("N5").Value = 1680
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 91
$ws.Range("Q5").Value = 43
$ws.Range("R5").Value = -198
$ws.Range("S5").Value = 53
$ws.Range("T5").Value = 147
$ws.Range("U5").Value = -104
$ws.Range("V5").Value = 100
$ws.Range("W5").Value = 9.529999999999999
$ws.Range("X5").Value = 3.87
$ws.Range("Y5").Value = 6.52
$ws.Range("Z5").Value = 4.65
$ws.Range("AA5").Value = 35.13
$ws.Range("AB5").Value = 1784.87
$ws.Range("AC5").Value = 546
$ws.Range("AD5").Value = 36.15
$ws.Range("AE5").Value = 8673
$ws.Range("AF5").Value = 2.28
$ws.Range("AG5").Value = 241
$ws.Range("AH5").Value = 1.22
$ws.Range("AI5").Value = 43.35
$ws.Range("AJ5").Value = 19662677

# Row 6
$ws.Range("D6").Value = 2867
$ws.Range("E6").Value = 307
$ws.Range("F6").Value = 307
$ws.Range("G6").Value = 312
$ws.Range("H6").Value = 233
$ws.Range("I6").Value = 237
$ws.Range("K6").Value = 2754
$ws.Range("L6").Value = 922
$ws.Range("M6").Value = 1832
$ws.Range("N6").Value = 1835
$ws.Range("P6").Value = 95
$ws.Range("Q6").Value = 400
$ws.Range("R6").Value = -417
$ws.Range("S6").Value = 129
$ws.Range("T6").Value = 408
$ws.Range("U6").Value = -8
$ws.Range("V6").Value = 280
$ws.Range("W6").Value = 10.73
$ws.Range("X6").Value = 8.119999999999999
$ws.Range("Y6").Value = 13.51
$ws.Range("Z6").Value = 9.27
$ws.Range("AA6").Value = 50.35
$ws.Range("AB6").Value = 1860.77
$ws.Range("AC6").Value = 1208
$ws.Range("AD6").Value = 12.99
$ws.Range("AE6").Value = 9476
$ws.Range("AF6").Value = 1.66
$ws.Range("AG6").Value = 253
$ws.Range("AH6").Value = 1.61
$ws.Range("AI6").Value = 20.59
$ws.Range("AJ6").Value = 19662677

# Row 7
$ws.Range("D7").Value = 3152
$ws.Range("E7").Value = 347
$ws.Range("G7").Value = 351
$ws.Range("H7").Value = 260
$ws.Range("I7").Value = 264
$ws.Range("K7").Value = 3187
$ws.Range("L7").Value = 1144
$ws.Range("M7").Value = 2044
$ws.Range("N7").Value = 2049
$ws.Range("P7").Value = 98
$ws.Range("Q7").Value = 313
$ws.Range("R7").Value = -415
$ws.Range("S7").Value = 107
$ws.Range("T7").Value = 303
$ws.Range("U7").Value = -52
$ws.Range("W7").Value = 11
$ws.Range("X7").Value = 8.26
$ws.Range("Y7").Value = 13.61
$ws.Range("Z7").Value = 8.76
$ws.Range("AA7").Value = 55.96
$ws.Range("AC7").Value = 1344
$ws.Range("AD7").Value = 11.79
$ws.Range("AE7").Value = 10581
$ws.Range("AF7").Value = 1.5
$ws.Range("AG7").Value = 267
$ws.Range("AH7").Value = 1.68
$ws.Range("AI7").Value = 19.84

# Row 8
$ws.Range("D8").Value = 3479
$ws.Range("E8").Value = 381
$ws.Range("G8").Value = 380
$ws.Range("H8").Value = 285
$ws.Range("I8").Value = 289
$ws.Range("K8").Value = 3424
$ws.Range("L8").Value = 1145
$ws.Range("M8").Value = 2278
$ws.Range("N8").Value = 2284
$ws.Range("P8").Value = 99
$ws.Range("Q8").Value = 434
$ws.Range("R8").Value = -197
$ws.Range("S8").Value = -91
$ws.Range("T8").Value = 133
$ws.Range("U8").Value = 194
$ws.Range("W8").Value = 10.96
$ws.Range("X8").Value = 8.19
$ws.Range("Y8").Value = 13.34
$ws.Range("Z8").Value = 8.619999999999999
$ws.Range("AA8").Value = 50.27
$ws.Range("AC8").Value = 1470
$ws.Range("AD8").Value = 10.78
$ws.Range("AE8").Value = 11797
$ws.Range("AF8").Value = 1.34
$ws.Range("AG8").Value = 287
$ws.Range("AH8").Value = 1.81
$ws.Range("AI8").Value = 19.5

# Row 9
$ws.Range("D9").Value = 3792
$ws.Range("E9").Value = 425
$ws.Range("G9").Value = 425
$ws.Range("H9").Value = 319
$ws.Range("I9").Value = 323
$ws.Range("K9").Value = 3676
$ws.Range("L9").Value = 1134
$ws.Range("M9").Value = 2542
$ws.Range("N9").Value = 2548
$ws.Range("P9").Value = 99
$ws.Range("Q9").Value = 414
$ws.Range("R9").Value = -210
$ws.Range("S9").Value = -110
$ws.Range("T9").Value = 145
$ws.Range("U9").Value = 214
$ws.Range("W9").Value = 11.21
$ws.Range("X9").Value = 8.4
$ws.Range("Y9").Value = 13.35
$ws.Range("Z9").Value = 8.98
$ws.Range("AA9").Value = 44.6
$ws.Range("AC9").Value = 1641
$ws.Range("AD9").Value = 9.66
$ws.Range("AE9").Value = 13158
$ws.Range("AF9").Value = 1.2
$ws.Range("AG9").Value = 293
$ws.Range("AH9").Value = 1.85
$ws.Range("AI9").Value = 17.88
